# Trade #25 closed at 2026-02-17 23:57:40 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets to reflect the newly closed trade #25 (MarketMaking strategy).

$wb = $excel.ActiveWorkbook

# Helper: write a string value into a cell without letting Excel's
# COM layer auto-convert date/time-looking text into a real date/time
# serial number (and without leaving a stray NumberFormat behind).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1500.82
$wsSummary.Range("B4").Value = 0.82
$wsSummary.Range("B5").Value = 0.66
$wsSummary.Range("B6").Value = 25
$wsSummary.Range("B7").Value = 14
$wsSummary.Range("B9").Value = 56

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 100.82
$wsStatus.Range("D6").Value = 25
$wsStatus.Range("E6").Value = 0.82
$wsStatus.Range("F6").Value = 0.82
$wsStatus.Range("G6").Value = 56

# ---------------------------------------------------------------------------
# New trade row (#25) appended to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------------
$newRow = 26

function Add-TradeRow($ws) {
    $ws.Cells.Item($newRow, 1).Value = 25

    Set-TextValue $ws.Cells.Item($newRow, 2) "2026-02-17"
    Set-TextValue $ws.Cells.Item($newRow, 3) "23:57:33"

    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "UP"
    $ws.Cells.Item($newRow, 6).Value = 0.75
    $ws.Cells.Item($newRow, 7).Value = 0.8
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = 6.6667
    $ws.Cells.Item($newRow, 10).Value = 0.05
    $ws.Cells.Item($newRow, 11).Value = 100.82
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.14
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
